$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-21 02:58:44"
$wsZhCn.Range("G5").Value = "2016-01-21 02:59:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-21 02:58:57"
$wsDeDe.Range("G5").Value = "2016-01-21 02:59:53"
